# Add 7 new rows (380-386) of landscaping data below the existing data,
# extending the table from A1:U379 to A1:U386, matching the style/format
# of the last existing row (379) and keeping the F-column ABS() formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the new rows 380-386 by copying the formatting (incl. the date
# number format in column A) from the last existing data row (379).
$ws.Range("A379:T379").Copy($ws.Range("A380:T386"))

# Row 380
$ws.Range("A380").Value = 45841
$ws.Range("B380").Value = "Flowering"
$ws.Range("C380").Value = "Large"
$ws.Range("D380").Value = 62
$ws.Range("E380").Value = 85
$ws.Range("F380").Formula = "=ABS(D380-E380)"
$ws.Range("G380").Value = 0
$ws.Range("H380").Value = 0
$ws.Range("I380").Value = "No"
$ws.Range("J380").Value = 2
$ws.Range("K380").Value = "Bright"
$ws.Range("L380").Value = 8
$ws.Range("M380").Value = 0.51
$ws.Range("N380").Value = 64
$ws.Range("O380").Value = 29.95
$ws.Range("P380").Value = 16
$ws.Range("Q380").Value = 0.27
$ws.Range("R380").Value = 9.9
$ws.Range("S380").Value = 54
$ws.Range("T380").Value = 0

# Row 381
$ws.Range("A381").Value = 45841
$ws.Range("B381").Value = "Nonflowering"
$ws.Range("C381").Value = "Medium"
$ws.Range("D381").Value = 62
$ws.Range("E381").Value = 85
$ws.Range("F381").Formula = "=ABS(D381-E381)"
$ws.Range("G381").Value = 0
$ws.Range("H381").Value = 0
$ws.Range("I381").Value = "No"
$ws.Range("J381").Value = 3
$ws.Range("K381").Value = "Bright"
$ws.Range("L381").Value = 8
$ws.Range("M381").Value = 0.51
$ws.Range("N381").Value = 64
$ws.Range("O381").Value = 29.95
$ws.Range("P381").Value = 16
$ws.Range("Q381").Value = 0.27
$ws.Range("R381").Value = 9.9
$ws.Range("S381").Value = 54
$ws.Range("T381").Value = 0

# Row 382
$ws.Range("A382").Value = 45841
$ws.Range("B382").Value = "Nonflowering"
$ws.Range("C382").Value = "Small"
$ws.Range("D382").Value = 62
$ws.Range("E382").Value = 85
$ws.Range("F382").Formula = "=ABS(D382-E382)"
$ws.Range("G382").Value = 0
$ws.Range("H382").Value = 0
$ws.Range("I382").Value = "No"
$ws.Range("J382").Value = 3
$ws.Range("K382").Value = "Bright"
$ws.Range("L382").Value = 8
$ws.Range("M382").Value = 0.51
$ws.Range("N382").Value = 64
$ws.Range("O382").Value = 29.95
$ws.Range("P382").Value = 16
$ws.Range("Q382").Value = 0.27
$ws.Range("R382").Value = 9.9
$ws.Range("S382").Value = 54
$ws.Range("T382").Value = 0

# Row 383
$ws.Range("A383").Value = 45841
$ws.Range("B383").Value = "Nonflowering"
$ws.Range("C383").Value = "Medium"
$ws.Range("D383").Value = 62
$ws.Range("E383").Value = 85
$ws.Range("F383").Formula = "=ABS(D383-E383)"
$ws.Range("G383").Value = 0
$ws.Range("H383").Value = 0.1
$ws.Range("I383").Value = "No"
$ws.Range("J383").Value = 3
$ws.Range("K383").Value = "Neutral"
$ws.Range("L383").Value = 8
$ws.Range("M383").Value = 0.51
$ws.Range("N383").Value = 64
$ws.Range("O383").Value = 29.95
$ws.Range("P383").Value = 16
$ws.Range("Q383").Value = 0.27
$ws.Range("R383").Value = 9.9
$ws.Range("S383").Value = 54
$ws.Range("T383").Value = 0

# Row 384
$ws.Range("A384").Value = 45841
$ws.Range("B384").Value = "Nonflowering"
$ws.Range("C384").Value = "Medium"
$ws.Range("D384").Value = 62
$ws.Range("E384").Value = 85
$ws.Range("F384").Formula = "=ABS(D384-E384)"
$ws.Range("G384").Value = 0
$ws.Range("H384").Value = 0.1
$ws.Range("I384").Value = "No"
$ws.Range("J384").Value = 3
$ws.Range("K384").Value = "Neutral"
$ws.Range("L384").Value = 8
$ws.Range("M384").Value = 0.51
$ws.Range("N384").Value = 64
$ws.Range("O384").Value = 29.95
$ws.Range("P384").Value = 16
$ws.Range("Q384").Value = 0.27
$ws.Range("R384").Value = 9.9
$ws.Range("S384").Value = 54
$ws.Range("T384").Value = 0

# Row 385
$ws.Range("A385").Value = 45841
$ws.Range("B385").Value = "Nonflowering"
$ws.Range("C385").Value = "Large"
$ws.Range("D385").Value = 62
$ws.Range("E385").Value = 85
$ws.Range("F385").Formula = "=ABS(D385-E385)"
$ws.Range("G385").Value = 0
$ws.Range("H385").Value = 0.2
$ws.Range("I385").Value = "No"
$ws.Range("J385").Value = 4
$ws.Range("K385").Value = "Bright"
$ws.Range("L385").Value = 8
$ws.Range("M385").Value = 0.51
$ws.Range("N385").Value = 64
$ws.Range("O385").Value = 29.95
$ws.Range("P385").Value = 16
$ws.Range("Q385").Value = 0.27
$ws.Range("R385").Value = 9.9
$ws.Range("S385").Value = 54
$ws.Range("T385").Value = 0

# Row 386
$ws.Range("A386").Value = 45841
$ws.Range("B386").Value = "Tree"
$ws.Range("C386").Value = "Medium"
$ws.Range("D386").Value = 62
$ws.Range("E386").Value = 85
$ws.Range("F386").Formula = "=ABS(D386-E386)"
$ws.Range("G386").Value = 0
$ws.Range("H386").Value = 0.3
$ws.Range("I386").Value = "No"
$ws.Range("J386").Value = 1
$ws.Range("K386").Value = "Neutral"
$ws.Range("L386").Value = 8
$ws.Range("M386").Value = 0.51
$ws.Range("N386").Value = 64
$ws.Range("O386").Value = 29.95
$ws.Range("P386").Value = 16
$ws.Range("Q386").Value = 0.27
$ws.Range("R386").Value = 9.9
$ws.Range("S386").Value = 54
$ws.Range("T386").Value = 0

# Update the view: scroll so row 379 is at the top and select N380:N386,
# matching where the author was working after pasting the new rows.
$excel.ActiveWindow.ScrollRow = 379
[void]$ws.Range("N380:N386").Select()
